# "Var är vi nu?" -> "Lägesrapport utvecklingen"
#
# The title placeholder text is updated on both slides that carried the
# old "Var är vi nu?" status-update heading (slides 7 and 8 in
# presentation order), matching the canonical OOXML change.
#
# (A reviewer "Oscar" modern/threaded comment was also attached to slide
# 7 in the original edit. This COM-interop runtime only creates legacy
# PowerPoint comments - not the modern `p188:cmLst` threaded-comment
# parts/relationships the real file uses - so faking it here would only
# add spurious, incorrectly-shaped parts instead of reproducing the
# change; the text edits below are the reproducible part of the diff.)

$p = $ppt.ActivePresentation

$newTitle = "Lägesrapport utvecklingen"

$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = $newTitle

$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = $newTitle
